# Automatische test-sync: 2025-06-17 21:37:03
# Append the new inbound mail log entry to "Logs" and refresh the
# "Dashboard" category tally accordingly.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 27 -------------------------------------------
$wsLogs = $wb.Worksheets.Item("Logs")

$wsLogs.Cells.Item(27, 1).Value = "Offerte voor zakelijke samenwerking"
$wsLogs.Cells.Item(27, 2).Value = "mailmind.test@zohomail.eu"
$wsLogs.Cells.Item(27, 3).Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$wsLogs.Cells.Item(27, 4).Value = "Bestelling"
$wsLogs.Cells.Item(27, 6).Value = "2025-06-17 21:36:58"
$wsLogs.Cells.Item(27, 7).Value = "Nee"

# Extend the conditional-formatting ranges (D2:D26 -> D2:D27, G2:G26 -> G2:G27)
# so the newly appended row keeps getting highlighted like the rest.
$wsLogs.Range("D2:D26").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D27"))
$wsLogs.Range("G2:G26").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G27"))

# --- Dashboard sheet: update category counts ------------------------------
# New row is category "Bestelling", so Bestelling now counts 2 (same as
# Klacht). Row 5 / Row 6 swap places: Bestelling moves to row 5, Klacht to
# row 6, both with a count of 2.
$wsDash = $wb.Worksheets.Item("Dashboard")

$wsDash.Cells.Item(5, 1).Value = "Bestelling"
$wsDash.Cells.Item(6, 1).Value = "Klacht"
$wsDash.Cells.Item(6, 2).Value = 2
